$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain assignments: safe for the "Volume(1h)" column (always text because of the
# "%"/spaces) and for the handful of "Price" values that contain two dots (e.g.
# "25.954.76") which Excel cannot parse as a number either way.
$ws.Range("D2").Value = "25.954.76"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.641.84"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "1.652.28"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").Value = "0.0₅7845"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "26.024.31"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").Value = "1.132.85"
$ws.Range("E37").Value = "  -2.30%  "
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("E39").Value = "  +15.05%  "
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").Value = "1.781.10"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("E51").Value = "  +0.51%  "

# The remaining "Price" values look like plain decimal numbers (e.g. "215.39").
# A direct .Value assignment would have Excel auto-convert them into a numeric
# cell, which does not match the original (they are stored as literal text in
# this sheet). Build each one as a text-formula result in a scratch cell, then
# copy/paste-special just the value into the target cell so it lands as text
# without disturbing the target cells style.
$helper = $ws.Range("ZZ1")

$helper.Formula = '="215.39"'
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$helper.Formula = '="0.5083"'
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$helper.Formula = '="1.004"'
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$helper.Formula = '="0.2564"'
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$helper.Formula = '="0.06383"'
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$helper.Formula = '="19.52"'
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$helper.Formula = '="0.07777"'
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$helper.Formula = '="4.300"'
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$helper.Formula = '="0.5467"'
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$helper.Formula = '="64.50"'
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$helper.Formula = '="1.004"'
$helper.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$helper.Formula = '="198.06"'
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$helper.Formula = '="4.443"'
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$helper.Formula = '="9.974"'
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$helper.Formula = '="6.068"'
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$helper.Formula = '="1.005"'
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$helper.Formula = '="1.880"'
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$helper.Formula = '="141.06"'
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$helper.Formula = '="0.1143"'
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$helper.Formula = '="6.885"'
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$helper.Formula = '="1.240"'
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$helper.Formula = '="0.05036"'
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$helper.Formula = '="3.265"'
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$helper.Formula = '="3.192"'
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$helper.Formula = '="1.545"'
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$helper.Formula = '="2.366"'
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$helper.Formula = '="0.8974"'
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$helper.Formula = '="2.597"'
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$helper.Formula = '="0.5500"'
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$helper.Formula = '="0.01559"'
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$helper.Formula = '="1.004"'
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$helper.Formula = '="5.622"'
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$helper.Formula = '="0.8212"'
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$helper.Formula = '="100.22"'
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$helper.Formula = '="0.4528"'
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$helper.Formula = '="1.005"'
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$helper.Formula = '="55.00"'
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$helper.Formula = '="0.05075"'
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)

$helper.Clear()
$excel.CutCopyMode = 0

